$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for column Q ---
$ws.Range("Q1").Value = "Running Time"

# --- Fill in the previously-blank prediction cells on row 38 ---
$ws.Range("M38").Value = 168.9
$ws.Range("N38").Value = 22.8
$ws.Range("O38").Value = 517
$ws.Range("P38").Value = 642.1370561497326

# --- New data row 39 ---
$ws.Range("A39").Value = "2017.06.07 23.49.04"
$ws.Range("B39").Value = 17
$ws.Range("C39").Value = 82
$ws.Range("D39").Value = 1015
$ws.Range("E39").Value = 3.099999904632568
$ws.Range("F39").Value = "2017-06-08T06:00:00"
$ws.Range("G39").Value = "2017-06-08T09:00:00"
$ws.Range("H39").Value = 14.77000045776367
$ws.Range("I39").Value = 997.6500244140625
$ws.Range("J39").Value = 94
$ws.Range("K39").Value = 0.004999999888241291
$ws.Range("L39").Value = 0.9200000166893005
$ws.Range("M39").Value = 163.7
$ws.Range("N39").Value = 24
$ws.Range("O39").Value = 458
$ws.Range("P39").Value = 782.1946466768526
$ws.Range("Q39").Value = 26.51946466768526

# --- New data row 40 ---
$ws.Range("A40").Value = "2017.06.08 01.19.29"
$ws.Range("B40").Value = 16.639999389648438
$ws.Range("C40").Value = 82
$ws.Range("D40").Value = 1016
$ws.Range("E40").Value = 1.5
$ws.Range("F40").Value = "2017-06-08T06:00:00"
$ws.Range("G40").Value = "2017-06-08T09:00:00"
$ws.Range("H40").Value = 14.479999542236328
$ws.Range("I40").Value = 997.6500244140625
$ws.Range("J40").Value = 94
$ws.Range("K40").Value = 0.004999999888241291
$ws.Range("L40").Value = 0.9200000166893005
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = ""
$ws.Range("O40").Value = ""

# --- Update the active selection to match where work left off ---
$ws.Range("N8").Select()
